# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" sheet and
#    populate it with the per-fund holdings detail for the new quarter.
# 2) Insert a new summary row for "2022-Q1" at the top of the "总计" sheet's
#    data (row 2), shifting the existing quarters down, and renumber the
#    leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q1" sheet, positioned right before "总计"
# ---------------------------------------------------------------------
$total0 = $wb.Worksheets.Item("总计")
$newSheet0 = $wb.Worksheets.Add($total0)
$newSheet0.Name = "2022-Q1"

# References can go stale across structural operations (sheet add/insert) -
# re-fetch by name after any such change before relying on them again.
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------
# Step 2: "2022-Q1" header row - reuse the bold/bordered header style
# already used by the "总计" sheet's header (column B onward, no A1).
# ---------------------------------------------------------------------
$total.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# Step 3: "2022-Q1" data rows. Column A is the bold/bordered running
# index (0-based); columns B-G are text (fund code/name/amounts kept as
# strings, matching the source data - e.g. "005994" keeps its leading
# zero and "11.53" stays textual); column H is a plain number (rank).
# ---------------------------------------------------------------------
$total.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

$fundRows = @(
    @("005994", "国投瑞银中证500指数量化增强A", "11.53", "87.00", "1.50", "0.1730", 3),
    @("007089", "国投瑞银中证500指数量化增强C", "3.82",  "87.00", "1.50", "0.0573", 3),
    @("011731", "国投瑞银安睿混合A",            "2.58",  "43.48", "0.74", "0.0191", 3),
    @("011732", "国投瑞银安睿混合C",            "0.95",  "43.48", "0.74", "0.0070", 3)
)

# Force column B-G to be stored as text so numeric-looking strings (fund
# codes, "11.53", "87.00", ...) are not coerced into numbers.
$newSheet.Range("B2:G5").NumberFormat = "@"

$r = 2
$idx = 0
foreach ($fundRow in $fundRows) {
    $newSheet.Range("A" + $r).Value = $idx
    $newSheet.Range("B" + $r).Value = $fundRow[0]
    $newSheet.Range("C" + $r).Value = $fundRow[1]
    $newSheet.Range("D" + $r).Value = $fundRow[2]
    $newSheet.Range("E" + $r).Value = $fundRow[3]
    $newSheet.Range("F" + $r).Value = $fundRow[4]
    $newSheet.Range("G" + $r).Value = $fundRow[5]
    $newSheet.Range("H" + $r).Value = $fundRow[6]
    $r++
    $idx++
}

# Drop the temporary "@" text number-format now that the values are safely
# stored as text, so the cells end up with plain/default formatting again
# (matches the source data, which carries no explicit style on B:G).
$newSheet.Range("B2:G5").ClearFormats()

# ---------------------------------------------------------------------
# Step 4: "总计" sheet - insert a new row 2 for "2022-Q1" (4 funds,
# 0.26 亿元), pushing the existing quarters down a row, then fix up the
# running index column (A) for the rows that moved.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.26

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
